$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 3
    3  = 1
    4  = 1
    5  = 0
    6  = 0
    7  = 2
    8  = 2
    9  = 0
    10 = 0
    11 = 1
    12 = 1
    13 = 0
    14 = 0
    15 = 1
    16 = 1
    17 = 2
    18 = 1
    19 = 0
    20 = 1
    21 = 0
    22 = 1
    23 = 1
    24 = 2
    25 = 0
    26 = 2
    27 = 3
    28 = 1
    29 = 3
    30 = 1
    31 = 1
    32 = 0
    33 = 5
    34 = 1
    36 = 2
    38 = 1
    39 = 1
    41 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
